# Add team record (Wins / Losses / Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 should look like the other headers (bold,
# centered, bordered) -- copy the formatting from the existing header A1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-41) gets the same team record: 84 wins, 78 losses, 0 ties.
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 30).Value = 84
    $ws.Cells.Item($r, 31).Value = 78
    $ws.Cells.Item($r, 32).Value = 0
}
